$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '74.939.07'
$ws.Range("E2").Value = '  +0.71%  '

$ws.Range("D3").Value = '2.810.33'
$ws.Range("E3").Value = '  +6.57%  '

$ws.Range("E4").Value = '  -0.04%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '187.15'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.04%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '594.94'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +1.71%  '

$ws.Range("E7").Value = '  -0.02%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.548'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +2.62%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.192'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -5.13%  '

$ws.Range("D10").Value = '2.810.09'
$ws.Range("E10").Value = '  +6.49%  '

$ws.Range("B11").Value = 'Cardano'
$ws.Range("C11").Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.374'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +4.31%  '

$ws.Range("B12").Value = 'TRON'
$ws.Range("C12").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.160'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -1.89%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.90'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +3.85%  '

$ws.Range("D14").Value = '3.328.32'
$ws.Range("E14").Value = '  +6.32%  '

$ws.Range("D15").Value = '74.914.78'
$ws.Range("E15").Value = '  +0.77%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0000187'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -1.77%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '26.87'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +1.66%  '

$ws.Range("D18").Value = '2.813.18'
$ws.Range("E18").Value = '  +6.82%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '9.03'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.66%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '12.28'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +3.86%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '377.96'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +1.47%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '2.28'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -1.47%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.09'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.07%  '

$ws.Range("E24").Value = '  +0.03%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '70.90'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.19%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.93'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +6.19%  '

$ws.Range("D27").Value = '2.955.89'
$ws.Range("E27").Value = '  +6.30%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '4.17'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.29%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.0000104'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +9.17%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.00'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.02%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '517.64'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -2.71%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.39'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.20%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '7.67'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.68%  '

$ws.Range("E34").Value = '  +1.33%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.999'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.05%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '163.19'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.21%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '19.95'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +3.91%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.119'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.02%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '19.39'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.49%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '185.75'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +15.62%  '

$ws.Range("E41").Value = '  +0.05%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.341'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +4.15%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.00'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +1.31%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.67'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.86%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.22'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +2.62%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '39.99'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +2.42%  '

$ws.Range("E47").Value = '  +0.47%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.33'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -2.80%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.580'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +9.24%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '3.72'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +2.85%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.635'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +8.22%  '
